$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Insert 3 new rows at 862 (pushes existing rows 862-920 down to 865-923)
$ws.Rows("862:864").Insert()

# Row 862: new record, Fecha=44610 (2022-02-18)
$ws.Cells.Item(862,1).Value = 3
$ws.Cells.Item(862,2).Value = "Femacal de La Calera"
$ws.Cells.Item(862,3).Value = "Coquimbo"
$ws.Cells.Item(862,4).Value = 44610
$ws.Cells.Item(862,5).Value = 5
$ws.Cells.Item(862,6).Value = 100112020
$ws.Cells.Item(862,7).Value = "Tomate"
$ws.Cells.Item(862,8).Value = "Larga vida"
$ws.Cells.Item(862,9).Value = "Extra"
$ws.Cells.Item(862,10).Value = 230
$ws.Cells.Item(862,11).Value = 16000
$ws.Cells.Item(862,12).Value = 17000
$ws.Cells.Item(862,13).Value = 16522
$ws.Cells.Item(862,14).Value = "$/bandeja 18 kilos"
$ws.Cells.Item(862,15).Value = "Limache"
$ws.Cells.Item(862,16).Value = 918
$ws.Cells.Item(862,17).Value = 18
$ws.Cells.Item(862,18).Value = "Hortaliza"

# Row 863: new record, Fecha=44610 (2022-02-18)
$ws.Cells.Item(863,1).Value = 3
$ws.Cells.Item(863,2).Value = "Femacal de La Calera"
$ws.Cells.Item(863,3).Value = "Coquimbo"
$ws.Cells.Item(863,4).Value = 44610
$ws.Cells.Item(863,5).Value = 5
$ws.Cells.Item(863,6).Value = 100112020
$ws.Cells.Item(863,7).Value = "Tomate"
$ws.Cells.Item(863,8).Value = "Larga vida"
$ws.Cells.Item(863,9).Value = "Primera"
$ws.Cells.Item(863,10).Value = 120
$ws.Cells.Item(863,11).Value = 14000
$ws.Cells.Item(863,12).Value = 14000
$ws.Cells.Item(863,13).Value = 14000
$ws.Cells.Item(863,14).Value = "$/bandeja 18 kilos"
$ws.Cells.Item(863,15).Value = "Limache"
$ws.Cells.Item(863,16).Value = 778
$ws.Cells.Item(863,17).Value = 18
$ws.Cells.Item(863,18).Value = "Hortaliza"

# Row 864: new record, Fecha=44610 (2022-02-18)
$ws.Cells.Item(864,1).Value = 3
$ws.Cells.Item(864,2).Value = "Femacal de La Calera"
$ws.Cells.Item(864,3).Value = "Coquimbo"
$ws.Cells.Item(864,4).Value = 44610
$ws.Cells.Item(864,5).Value = 5
$ws.Cells.Item(864,6).Value = 100112020
$ws.Cells.Item(864,7).Value = "Tomate"
$ws.Cells.Item(864,8).Value = "Larga vida"
$ws.Cells.Item(864,9).Value = "Segunda"
$ws.Cells.Item(864,10).Value = 230
$ws.Cells.Item(864,11).Value = 12000
$ws.Cells.Item(864,12).Value = 12500
$ws.Cells.Item(864,13).Value = 12261
$ws.Cells.Item(864,14).Value = "$/bandeja 18 kilos"
$ws.Cells.Item(864,15).Value = "Limache"
$ws.Cells.Item(864,16).Value = 681
$ws.Cells.Item(864,17).Value = 18
$ws.Cells.Item(864,18).Value = "Hortaliza"
